# --- Refined metadata to be additional tab ---
$wb = $excel.ActiveWorkbook
$ds = $wb.Worksheets.Item("data")

# 1. Refresh the "time_taken" timestamps in the data sheet (re-run recorded
#    at 2021-10-05 14:35:07.7xxxxx instead of the original 10:52:06.0xxxxx).
$newTimes = @(
  "2021-10-05 14:35:07.771239",
  "2021-10-05 14:35:07.771247",
  "2021-10-05 14:35:07.771251",
  "2021-10-05 14:35:07.771254",
  "2021-10-05 14:35:07.771257",
  "2021-10-05 14:35:07.771259",
  "2021-10-05 14:35:07.771262",
  "2021-10-05 14:35:07.771264",
  "2021-10-05 14:35:07.771268",
  "2021-10-05 14:35:07.771270",
  "2021-10-05 14:35:07.771273",
  "2021-10-05 14:35:07.771275",
  "2021-10-05 14:35:07.771278",
  "2021-10-05 14:35:07.771281",
  "2021-10-05 14:35:07.771283",
  "2021-10-05 14:35:07.771286",
  "2021-10-05 14:35:07.771289",
  "2021-10-05 14:35:07.771292",
  "2021-10-05 14:35:07.771294",
  "2021-10-05 14:35:07.771297",
  "2021-10-05 14:35:07.771300",
  "2021-10-05 14:35:07.771302",
  "2021-10-05 14:35:07.771305",
  "2021-10-05 14:35:07.771307",
  "2021-10-05 14:35:07.771310",
  "2021-10-05 14:35:07.771313",
  "2021-10-05 14:35:07.771316",
  "2021-10-05 14:35:07.771318",
  "2021-10-05 14:35:07.771321",
  "2021-10-05 14:35:07.771324",
  "2021-10-05 14:35:07.771326",
  "2021-10-05 14:35:07.771329",
  "2021-10-05 14:35:07.771332",
  "2021-10-05 14:35:07.771335",
  "2021-10-05 14:35:07.771337",
  "2021-10-05 14:35:07.771340",
  "2021-10-05 14:35:07.771343",
  "2021-10-05 14:35:07.771345",
  "2021-10-05 14:35:07.771348",
  "2021-10-05 14:35:07.771350"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
  $row = $i + 2
  $ds.Cells.Item($row, 6).Value = $newTimes[$i]
}

# 2. Add a new "metadata" worksheet after "data".
$ws = $wb.Worksheets.Add($null, $ds)
$ws.Name = "metadata"

# Header row.
$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

# Match header/border style used on the data tab's header row.
$ds.Range("B1:F1").Copy()
$ws.Range("B1:F1").PasteSpecial(-4122)
$ds.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)

# Data row.
$ds.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A2").Value = 0

$ws.Range("B2").Value = "Osteogenesis Imperfecta"
$ws.Range("C2").Value = 147
# Force "0.65" to be stored as text (matches source data), then restore the
# cell's style to the sheet's default (plain, unstyled) by re-copying the
# format of the equivalent cell from the data tab, which carries no style.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "0.65"
$ds.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("E2").Value = "2021-09-13T07:51:02.808727Z"
$ws.Range("F2").Value = "2021-10-05 14:35:07.767627"
$ws.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/147/?format=json"
